$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.237.95"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.857.13"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'232.70"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4714"
$ws.Range("D8").Value = "'0.2734"
$ws.Range("E8").Value = "  -4.18%  "
$ws.Range("D9").Value = "'0.06420"
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").Value = "1.843.81"
$ws.Range("E10").Value = "  -4.90%  "
$ws.Range("D11").Value = "'0.07446"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").Value = "  -3.04%  "
$ws.Range("D13").Value = "'4.992"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").Value = "'85.35"
$ws.Range("E14").Value = "  -4.02%  "
$ws.Range("D15").Value = "'0.6306"
$ws.Range("E15").Value = "  -5.97%  "
$ws.Range("D16").Value = "30.191.81"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "'0.9996"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'232.24"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'12.72"
$ws.Range("E19").Value = "  -4.88%  "
$ws.Range("D20").Value = "'0.000007359"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "2.092.86"
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'5.027"
$ws.Range("D24").Value = "'5.993"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").Value = "'9.242"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").Value = "'164.82"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").Value = "'17.86"
$ws.Range("E27").Value = "  -4.93%  "
$ws.Range("D28").Value = "'1.885"
$ws.Range("E28").Value = "  -3.30%  "
$ws.Range("D29").Value = "'0.1025"
$ws.Range("E29").Value = "  +5.30%  "
$ws.Range("D30").Value = "'1.381"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("D31").Value = "'4.131"
$ws.Range("E31").Value = "  -5.60%  "
$ws.Range("D32").Value = "'3.923"
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("D33").Value = "'0.04894"
$ws.Range("E33").Value = "  -3.49%  "
$ws.Range("D34").Value = "'1.144"
$ws.Range("E34").Value = "  -5.85%  "
$ws.Range("D35").Value = "'0.7220"
$ws.Range("E35").Value = "  -4.35%  "
$ws.Range("D36").Value = "'0.9992"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'2.692"
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "'0.01898"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").Value = "'2.638"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'0.9054"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").Value = "'1.970"
$ws.Range("E41").Value = "  -5.97%  "
$ws.Range("D42").Value = "'105.58"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "'0.4100"
$ws.Range("E44").Value = "  -4.66%  "
$ws.Range("D45").Value = "'5.513"
$ws.Range("E45").Value = "  -5.18%  "
$ws.Range("D46").Value = "'7.111"
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("D47").Value = "'61.00"
$ws.Range("E47").Value = "  -6.19%  "
$ws.Range("D48").Value = "'0.1197"
$ws.Range("E48").Value = "  -7.10%  "
$ws.Range("D49").Value = "'8.694"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("D50").Value = "'1.406"
$ws.Range("E50").Value = "  -5.36%  "
$ws.Range("D51").Value = "'33.24"
$ws.Range("E51").Value = "  -2.05%  "
